$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-orders the weekly price records (rows 2-12, excluding the
# already-sorted rows 4 and 9) by swapping / rotating the D..T data
# across rows, as described by the diff. We capture the "before" values
# for every affected row/column first, then write the "after" values,
# so that reads always see original data regardless of write order.

$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value2
    }
    return $vals
}

# Capture original values of every row that participates in the move.
$orig = @{}
foreach ($r in 2,3,5,6,7,8,10,11,12) {
    $orig[$r] = Get-RowValues $r
}

# Mapping: new content of row N comes from the original content of row M.
$mapping = @{
    2  = 10
    3  = 7
    5  = 3
    6  = 8
    7  = 6
    8  = 12
    10 = 2
    11 = 5
    12 = 11
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $orig[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcVals[$col]
    }
}
